# The presentation currently has two themes:
#   ppt/theme/theme1.xml -> used by the Slide Master, clrScheme name="Integral"
#   ppt/theme/theme2.xml -> used by the Notes Master, clrScheme name="Office"
# The target edit swaps the two themes' color schemes (and, in the canonical
# XML, their name attributes) so theme1.xml ends up holding the "Office"
# colors and theme2.xml ends up holding the "Integral" colors.
#
# Colors are set through ThemeColorScheme.Colors(i).RGB (the documented,
# supported way to edit a theme's palette in this host). The VBA/OLE RGB
# encoding packs R in the low byte, G in the middle byte and B in the high
# byte, i.e. RGB(r,g,b) = r + g*256 + b*65536 -- the reverse of the usual
# "RRGGBB" hex string order.

$p = $ppt.ActivePresentation

$m  = $p.SlideMaster
$nm = $p.NotesMaster

$mScheme  = $m.Theme.ThemeColorScheme
$nmScheme = $nm.Theme.ThemeColorScheme

# Target palette for the Slide Master (was "Integral", becomes "Office").
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

# Target palette for the Notes Master (was "Office", becomes "Integral").
$integralColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x515F45,  # dk2      455F51
    0xD1DEE3,  # lt2      E3DED1
    0x38CB99,  # accent1  99CB38
    0x37A563,  # accent2  63A537
    0x24D0E6,  # accent3  E6D024
    0x0097CC,  # accent4  CC9700
    0xCFB34E,  # accent5  4EB3CF
    0xA68D37,  # accent6  378DA6
    0x259F6B,  # hlink    6B9F25
    0x026BB2   # folHlink B26B02
)

for ($i = 1; $i -le 12; $i++) {
    $mScheme.Colors($i).RGB  = $officeColors[$i - 1]
    $nmScheme.Colors($i).RGB = $integralColors[$i - 1]
}
